# Trade #11 closed at 2026-02-17 07:58:27 - MarketMaking UP -33.333%
#
# Updates the "Summary" and "Strategy Status" roll-up figures, and appends
# the new trade row (Trade # 11) to both the "All Trades" and "MarketMaking"
# logs.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.96                 # Current Capital
$summary.Range("B4").Value = -0.04                   # Total P&L $
$summary.Range("B5").Value = -0.07000000000000001    # Total P&L %
$summary.Range("B6").Value = 11                      # Total Trades
$summary.Range("B8").Value = 6                       # Losing Trades
$summary.Range("B9").Value = 36.36                   # Win Rate %

# --- Strategy Status sheet (MarketMaking row, row 4) ------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95999999999999        # Capital
$status.Range("D4").Value = 11                       # Trades
$status.Range("E4").Value = -0.04                    # P&L $
$status.Range("F4").Value = -0.04                    # P&L %
$status.Range("G4").Value = 36.36                    # Win Rate %

# --- Append new trade row (#11) to "All Trades" and "MarketMaking" ----
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 12

    $ws.Cells.Item($row, 1).Value = 11

    # Force the date to stay literal text (matches how the other rows in
    # this log were written) instead of being auto-converted to a date
    # serial number; reset the style afterwards so no extra formatting
    # lingers on the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "07:58:27"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.09
    $ws.Cells.Item($row, 7).Value = 0.06
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -33.3333
    $ws.Cells.Item($row, 10).Value = -0.03
    $ws.Cells.Item($row, 11).Value = 99.95999999999999
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
